# Move regression estimates to country specific folder
# Rename the "EL" worksheet to "Population_projections"

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("EL")
$ws.Name = "Population_projections"
